# diagrams/data_flow_diagram.pptx - "inputs and flow diagram"
# Lower-cases several of the diagram's box labels.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Output Folder" -> "output"
# (TextBox 3 also contains a bulleted list below the title line; only the
# first run of the first paragraph - the "Output Folder" title - changes.)
$grp1 = $s.Shapes.Item("Group 6")
$tb1 = $grp1.GroupItems.Item("TextBox 3")
$tb1.TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1).Text = "output"

# "Objective Variables" -> "objective variables"
$grp2 = $s.Shapes.Item("Group 8")
$tb2 = $grp2.GroupItems.Item("TextBox 10")
$tb2.TextFrame.TextRange.Text = "objective variables"

# "Player Inputs" -> "player inputs"
$grp3 = $s.Shapes.Item("Group 11")
$tb3 = $grp3.GroupItems.Item("TextBox 13")
$tb3.TextFrame.TextRange.Text = "player inputs"

# "Data Transformation" -> "data transformation"
$grp4 = $s.Shapes.Item("Group 15")
$tb4 = $grp4.GroupItems.Item("TextBox 17")
$tb4.TextFrame.TextRange.Text = "data transformation"

# "Draft Inputs" -> "draft inputs"
$grp5 = $s.Shapes.Item("Group 79")
$tb5 = $grp5.GroupItems.Item("TextBox 81")
$tb5.TextFrame.TextRange.Text = "draft inputs"
